$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q3" sheet right before the existing "2022-Q2" sheet.
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

# Header row (row 1, columns B..H) - copy the bold/bordered header style from
# the neighbouring "2022-Q2" sheet so the new sheet matches the existing look.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$existingQ2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Data rows: A (index, number), B..G (text-like codes/percentages kept as
# text so leading zeros / trailing zeros survive), H (rank, number).
$rows = @(
        @(0, "012159", "财通资管健康产业混合A", "10.00", "94.52", "6.18", "0.6180", 6),
        @(1, "001559", "天弘医疗健康混合C", "4.16", "84.58", "3.37", "0.1402", 10),
        @(2, "002300", "长盛医疗行业量化配置股票", "2.39", "93.48", "5.77", "0.1379", 5),
        @(3, "012160", "财通资管健康产业混合C", "1.95", "94.52", "6.18", "0.1205", 6),
        @(4, "001558", "天弘医疗健康混合A", "2.84", "84.58", "3.37", "0.0957", 10),
        @(5, "000684", "长盛养老健康产业灵活配置混合", "1.35", "92.15", "5.52", "0.0745", 4),
        @(6, "014126", "华夏中证1000指数增强C", "8.78", "89.62", "0.79", "0.0694", 9),
        @(7, "008412", "长盛竞争优势股票A", "0.72", "91.35", "5.54", "0.0399", 4),
        @(8, "010434", "红土创新医疗保健股票", "0.44", "94.49", "8.60", "0.0378", 1),
        @(9, "008619", "永赢医药健康股票C", "0.40", "94.40", "8.27", "0.0331", 3),
        @(10, "008413", "长盛竞争优势股票C", "0.39", "91.35", "5.54", "0.0216", 4),
        @(11, "008618", "永赢医药健康股票A", "0.24", "94.40", "8.27", "0.0198", 3),
        @(12, "015139", "泰康医疗健康股票A", "0.34", "85.50", "2.69", "0.0091", 9),
        @(13, "014125", "华夏中证1000指数增强A", "0.97", "89.62", "0.79", "0.0077", 9),
        @(14, "008884", "博远博锐混合A", "0.13", "75.55", "4.38", "0.0057", 1),
        @(15, "015140", "泰康医疗健康股票C", "0.13", "85.50", "2.69", "0.0035", 9),
        @(16, "008885", "博远博锐混合C", "0.02", "75.55", "4.38", "0.0009", 1),
    )

# Columns B:G hold text that looks numeric ("012159", "10.00", ...) - format
# as Text first so assigning the value does not silently coerce to a number.
$q3.Range("B2:G18").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Column A (index numbers) uses the same bold/bordered "index" style as every
# other sheet - copy it from the neighbouring sheet's A2 cell.
$existingQ2.Range("A2").Copy()
$q3.Range("A2:A18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new top row for 2022-Q3 and
#    shift the existing 2022-Q2 / 2022-Q1 / 2021-Q4 rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A5").Value = 3
$total.Range("A2").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.01

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 14
$total.Range("D4").Value = 1.44

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 19
$total.Range("D3").Value = 1.96

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 17
$total.Range("D2").Value = 1.44
